$d = $word.ActiveDocument

$replacements = @(
    @{old = "522÷8="; new = "600÷2="},
    @{old = "862÷5="; new = "380÷8="},
    @{old = "820÷3="; new = "463÷3="},
    @{old = "497÷9="; new = "446÷8="},
    @{old = "866÷7="; new = "295÷7="},
    @{old = "196÷6="; new = "824÷4="},
    @{old = "484÷4="; new = "406÷3="},
    @{old = "340÷9="; new = "451÷2="},
    @{old = "592÷2="; new = "293÷5="},
    @{old = "660÷6="; new = "198÷6="},
    @{old = "699÷5="; new = "431÷2="},
    @{old = "897÷8="; new = "919÷8="},
    @{old = "965÷9="; new = "237÷5="},
    @{old = "248÷6="; new = "558÷3="},
    @{old = "170÷6="; new = "851÷9="},
    @{old = "507÷4="; new = "376÷6="},
    @{old = "427÷4="; new = "968÷6="},
    @{old = "294÷6="; new = "164÷4="},
    @{old = "281÷2="; new = "849÷3="},
    @{old = "286÷4="; new = "900÷8="},
    @{old = "365÷2="; new = "509÷7="},
    @{old = "661÷3="; new = "526÷4="},
    @{old = "120÷9="; new = "727÷8="},
    @{old = "245÷5="; new = "981÷6="},
    @{old = "964÷2="; new = "418÷7="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
